$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.672.57"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.115.40"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.59"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5257"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4540"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.02"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09085"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.170"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.38"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.122.82"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.809"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.082"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.72"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001163"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.015"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06721"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.420"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.755.98"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.380"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354.14"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.39"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.39"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.540"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.83"
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.200"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1077"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.367"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.632"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.951"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.938"
$ws.Range("E37").Value = "  +7.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02658"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06846"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2322"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.60"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.06"
$ws.Range("E44").Value = "  +6.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6434"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.312"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("E47").Value = "  +14.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.708"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.255"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07315"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.80"
$ws.Range("E51").Value = "  -0.76%  "
